# Adds a new "Layout" worksheet (architecture/markdown notes), positioned
# right after "Sheet3" and before "Sheet2", and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- Create the new sheet and put it in the right place ------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Layout"
$newSheet.Move($null, $wb.Worksheets.Item("Sheet3"))

# Re-acquire a fresh handle after the Move() (the old reference stops
# tracking the sheet's state correctly once it has been repositioned).
$layout = $wb.Worksheets.Item("Layout")

# --- Column widths (best achievable approximation of the authored sizes) -
$layout.Columns.Item(2).ColumnWidth = 20.5               # B ~ 21.29
$layout.Columns.Item(3).ColumnWidth = 11.5                # C ~ 12.29
$layout.Columns.Item(4).ColumnWidth = 12.6666666666667    # D ~ 13.57
$layout.Columns.Item(5).ColumnWidth = 12.8333333333333    # E ~ 13.71

# --- Cell content ------------------------------------------------------
# (written in the same order the original author typed it in, so freshly
# introduced shared strings land at the same indices as the authored file)
$layout.Range("B3").Value = "Server"
$layout.Range("M3").Value = "Client"

$layout.Range("B5").Value = "API for Fetch Data"

$layout.Range("B6").Value = "Web Application"
$layout.Range("B7").Value = "Management"
$layout.Range("B8").Value = "User creation"

$layout.Range("M4").Value = "Store Front"
$layout.Range("M5").Value = "Add Items, sell items"
$layout.Range("M6").Value = "Create and submit orders"
$layout.Range("M7").Value = "Close cash register"
$layout.Range("M8").Value = "create report related to cash register based on user level"
$layout.Range("M9").Value = "control user interface"
$layout.Range("M10").Value = "app settings like colors and language"

$layout.Range("F7").Value = "Stocks"
$layout.Range("G7").Value = "Printers"
$layout.Range("M11").Value = "setting up print stations"

$layout.Range("C7").Value = "Dashboard"
$layout.Range("D7").Value = "Documents"
$layout.Range("E7").Value = "Products"

$layout.Range("C8").Value = "Create and manage users"
$layout.Range("E8").Value = "Permissions and authorizations"

$layout.Range("C6").Value = "Like the store"
$layout.Range("E6").Value = "create and submit orders"

# --- View state: selection + active tab -------------------------------
$layout.Range("D11").Select()
$layout.Activate()
